# Auto-generated Excel COM-interop script
# Applies the 'update to include tokenization process' changes

$wb = $excel.ActiveWorkbook
$wsEntities = $wb.Worksheets.Item("Entities")
$wsColumns = $wb.Worksheets.Item("Columns")

# ---- Entities sheet: rewrite A1:B14 with the updated entity list ----
$entities = @(
    @("entity_name", "purpose"),
    @("TBL_VOID_UserSessionInput", "various inputs users put into Abyss"),
    @("TBL_SentimentAnalysis", "results of sentiment analysis (positive and negative) for user input"),
    @("TBL_SentimentAnalysis_Iteration", "configuration parameters for sentiment analysis iteration"),
    @("VW_SentimentRatings", "recombing sentiment analyis results with other metadata (text value, date, etc)"),
    @("SP_DeleteStagingTables_SentimentAnalysis", "delete relevant staging tables in sentiment analysis process"),
    @("SP_InsertInto_SentimentAnalysis_Iteration", "insert registry into sentiment analysis iteration table"),
    @("SP_InsertInto_SentimentAnalysis", "insert values into sentiment analysis table from staging table"),
    @("TBL_Tokens", "tokens extracted from tokenization process"),
    @("TBL_Tokens_Iteration", "configuration parameters for tokenization process"),
    @("VW_Tokens", "recombining tokenization results with other metadata (date, input text, etc)"),
    @("SP_DeleteStagingTables_SentimentAnalysis", "delete relevant staging tables in tokenization process"),
    @("SP_InsertInto_SentimentAnalysis_Iteration", "insert registry into tokenization iteration table"),
    @("SP_InsertInto_SentimentAnalysis", "insert values into tokenization table from staging table")
)

for ($r = 0; $r -lt $entities.Length; $r++) {
    $row = $entities[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsEntities.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---- Columns sheet: rewrite A1:J30 with the updated column catalogue ----
$columns = @(
    @("table_schema", "table_name", "column_name", "column_default", "is_nullable", "data_type", "primary_key", "foreign_key", "links_to_table", "links_to_column"),
    @("dbo", "TBL_VOID_UserSessionInput", "PK_ID_UserSessionInput", "NULL", "NO", "int", "yes", "no", "NULL", "no"),
    @("dbo", "TBL_VOID_UserSessionInput", "TS_UserSessionInput", "NULL", "YES", "datetime2", "no", "no", "NULL", "no"),
    @("dbo", "TBL_VOID_UserSessionInput", "ID_User", "NULL", "YES", "int", "no", "yes", "TBL_VOID_User", "PK_ID_User"),
    @("dbo", "TBL_VOID_UserSessionInput", "ID_Session", "NULL", "YES", "int", "no", "yes", "TBL_VOID_Session", "PK_ID_Session"),
    @("dbo", "TBL_VOID_UserSessionInput", "Input", "NULL", "YES", "nvarchar", "no", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis", "PK_ID_SentimentAnalysis", "NULL", "NO", "int", "yes", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis", "TS_SentimentAnalysis", "NULL", "YES", "datetime2", "no", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis", "FK_ID_SentimentAnalysis_Iteration", "NULL", "YES", "int", "no", "yes", "TBL_SentimentAnalysis_Iteration", "PK_ID_SentimentAnalysis_Iteration"),
    @("dbo", "TBL_SentimentAnalysis", "FK_ID_UserSessionInput", "NULL", "YES", "int", "no", "yes", "TBL_UserSessionInput", "PK_ID_UserSessionInput"),
    @("dbo", "TBL_SentimentAnalysis", "Positive_Sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis", "Negative_Sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("STG", "TBL_SentimentAnalysis", "FK_ID_UserSessionInput", "NULL", "YES", "int", "no", "yes", "TBL_UserSession_Input", "PK_ID_UserSessionInput"),
    @("STG", "TBL_SentimentAnalysis", "Positive_Sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("STG", "TBL_SentimentAnalysis", "Negative_Sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis_Iteration", "PK_ID_SentimentAnalysis_Iteration", "NULL", "NO", "int", "yes", "no", "NULL", "no"),
    @("dbo", "TBL_SentimentAnalysis_Iteration", "TS_SentimentAnalysis_Iteration", "NULL", "YES", "datetime2", "no", "no", "NULL", "no"),
    @("dbo", "VW_SentimentRatings", "date", "NULL", "YES", "datetime2", "no", "no", "NULL", "no"),
    @("dbo", "VW_SentimentRatings", "input", "NULL", "YES", "nvarchar", "no", "no", "NULL", "no"),
    @("dbo", "VW_SentimentRatings", "positive_sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("dbo", "VW_SentimentRatings", "negative_sentiment", "NULL", "YES", "float", "no", "no", "NULL", "no"),
    @("dbo", "TBL_Tokens", "PK_ID_Token", "NULL", "NO", "int", "yes", "no", "NULL", "no"),
    @("dbo", "TBL_Tokens", "TS_Token", "NULL", "YES", "datetime2", "no", "no", "NULL", "no"),
    @("dbo", "TBL_Tokens", "FK_ID_Token_Iteration", "NULL", "YES", "int", "no", "yes", "TBL_Tokens_Iteration", "PK_ID_Tokens_Iteration"),
    @("dbo", "TBL_Tokens", "FK_ID_UserSessionInput", "NULL", "YES", "int", "no", "yes", "TBL_VOID_UserSessionInput", "PK_ID_UserSessionInput"),
    @("dbo", "TBL_Tokens", "Token", "NULL", "YES", "nvarchar", "no", "no", "NULL", "no"),
    @("STG", "TBL_Tokens", "FK_ID_UserSessionInput", "NULL", "YES", "int", "no", "yes", "TBL_VOID_UserSessionInput", "PK_ID_UserSessionInput"),
    @("STG", "TBL_Tokens", "Token", "NULL", "YES", "nvarchar", "no", "no", "NULL", "no"),
    @("dbo", "TBL_Tokens_Iteration", "PK_ID_Tokens_Iteration", "NULL", "NO", "int", "yes", "no", "NULL", "no"),
    @("dbo", "TBL_Tokens_Iteration", "TS_Tokens_Iteration", "NULL", "YES", "datetime2", "no", "no", "NULL", "no")
)

for ($r = 0; $r -lt $columns.Length; $r++) {
    $row = $columns[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsColumns.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---- Selection / active-view bookkeeping to match the final workbook state ----
$wsColumns.Activate()
$wsColumns.Range("E19").Select()
$wsEntities.Range("B15").Select()
$wsColumns.Activate()

Write-Output "done"